$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (K2:T2) with new TPM-derived figures
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05619066666666667
$ws.Range("N2").Value = 0.168572
$ws.Range("O2").Value = 0.3931387525216601
$ws.Range("P2").Value = 0.39313875252166
$ws.Range("Q2").Value = 0.03133026747377778
$ws.Range("R2").Value = 0.281972407264
$ws.Range("S2").Value = 0.3931387525216601
$ws.Range("T2").Value = 0.39313875252166

# Update row 3 values (O3, P3, S3, T3) with new TPM-derived figures
$ws.Range("O3").Value = 0.60686124747834
$ws.Range("P3").Value = 0.60686124747834
$ws.Range("S3").Value = 0.60686124747834
$ws.Range("T3").Value = 0.60686124747834

# Remove row 4 entirely (the MuSCs target-cluster row no longer present)
$ws.Rows.Item(4).Delete()
